$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.697.26'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '3.705.30'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''677.58'
$ws.Range("D6").Value = '''161.83'
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("E9").Value = '  +1.71%  '
$ws.Range("D10").Value = '''7.14'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '''32.82'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").Value = '3.684.91'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = '69.722.14'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '''16.09'
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("D19").Value = '''473.88'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '''9.82'
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").Value = '''80.46'
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").Value = '3.852.79'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = '''0.0000127'
$ws.Range("E24").Value = '  +2.78%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '''11.01'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''26.97'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '''0.166'
$ws.Range("E34").Value = '  +4.27%  '
$ws.Range("D35").Value = '3.694.81'
$ws.Range("E35").Value = '  +0.84%  '
$ws.Range("E36").Value = '  +4.08%  '
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").Value = '''0.945'
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").Value = '''166.79'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").Value = '''28.31'
$ws.Range("E46").Value = '  +1.19%  '
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -1.46%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").Value = '''7.91'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("E51").Value = '  +2.01%  '
